$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that were deleted from the source data ---
# Row 26 = "RM 232", Row 28 (becomes 27 after first delete) = "SC 92"
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- Individual cell value corrections (rows 2-25, unaffected by the row deletions) ---
$ws.Range("F2").Value = ""
$ws.Range("F5").Value = 17.66
$ws.Range("C6").Value = 15.1
$ws.Range("F6").Value = 16.43
$ws.Range("C8").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("C12").Value = 12.5
$ws.Range("F13").Value = ""
$ws.Range("C14").Value = ""
$ws.Range("C17").Value = 11.2
$ws.Range("C18").Value = 11.5
$ws.Range("C19").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("C23").Value = 12.2
$ws.Range("F24").Value = 16.78

# --- Individual cell value corrections for rows 26-33 (post row-deletion numbering) ---
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = ""
$ws.Range("F28").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("F30").Value = 16.89
$ws.Range("B32").Value = ""
